# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45204 (2023-10-05) to 45205 (2023-10-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 271; $row++) {
    $ws.Cells.Item($row, 3).Value = 45205
}
